# Auto-generated edit script applying stock-qty/value corrections
# per the commit diff. Each data row's F (qty) and G (value = D*F)
# are corrected, and B (row/sub-total/grand-total amounts) updated
# to match. Some adjacent row pairs have their B/E/F/G values
# swapped (e.g. rows 192/193) rather than independently changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: F98->94, G2928.24->2808.72
$ws.Range("F6").Value = 94
$ws.Range("G6").Value = 2808.72

# Row 7: F96->95, G4495.68->4448.85
$ws.Range("F7").Value = 95
$ws.Range("G7").Value = 4448.85

# Row 10: B31106.79->30940.44
$ws.Range("B10").Value = 30940.44

# Row 59: F21->20, G1724.94->1642.8
$ws.Range("F59").Value = 20
$ws.Range("G59").Value = 1642.8

# Row 68: F59->58, G6792.08->6676.96
$ws.Range("F68").Value = 58
$ws.Range("G68").Value = 6676.96

# Row 77: F287->286, G13414.38->13367.64
$ws.Range("F77").Value = 286
$ws.Range("G77").Value = 13367.64

# Row 78: F45->44, G2560.5->2503.6
$ws.Range("F78").Value = 44
$ws.Range("G78").Value = 2503.6

# Row 80: F16->15, G3937.12->3691.05
$ws.Range("F80").Value = 15
$ws.Range("G80").Value = 3691.05

# Row 81: F28->25, G856.24->764.5
$ws.Range("F81").Value = 25
$ws.Range("G81").Value = 764.5

# Row 90: B201482.09->200843.38
$ws.Range("B90").Value = 200843.38

# Row 115: F230->229, G22266.3->22169.49
$ws.Range("F115").Value = 229
$ws.Range("G115").Value = 22169.49

# Row 117: B16318.58->16221.77
$ws.Range("B117").Value = 16221.77

# Row 144: F1224->1218, G10342.8->10292.1
$ws.Range("F144").Value = 1218
$ws.Range("G144").Value = 10292.1

# Row 145: F656->642, G5241.44->5129.58
$ws.Range("F145").Value = 642
$ws.Range("G145").Value = 5129.58

# Row 147: B18615.08->18452.52
$ws.Range("B147").Value = 18452.52

# Row 152: F73->71, G6445.17->6268.59
$ws.Range("F152").Value = 71
$ws.Range("G152").Value = 6268.59

# Row 156: B35636.7->35460.12
$ws.Range("B156").Value = 35460.12

# Row 192: B48706->64973, E39.8->35.4, F-144->2, G-4795.2->66.59999999999999
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.59999999999999

# Row 193: B64973->48706, E35.4->39.8, F2->-144, G66.59999999999999->-4795.2
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2

# Row 194: F44->43, G3771.68->3685.96
$ws.Range("F194").Value = 43
$ws.Range("G194").Value = 3685.96

# Row 205: F31->30, G11691.34->11314.2
$ws.Range("F205").Value = 30
$ws.Range("G205").Value = 11314.2

# Row 212: F10->9, G4671.9->4204.71
$ws.Range("F212").Value = 9
$ws.Range("G212").Value = 4204.71

# Row 216: B48943.91->48013.86
$ws.Range("B216").Value = 48013.86

# Row 222: F20->19, G2898.6->2753.67
$ws.Range("F222").Value = 19
$ws.Range("G222").Value = 2753.67

# Row 227: B55373->63520, E163.62->153.4, F-94->67, G-13562.32->9666.76
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 67
$ws.Range("G227").Value = 9666.76

# Row 228: B63520->55373, E153.4->163.62, F67->-94, G9666.76->-13562.32
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32

# Row 230: F67->66, G9613.16->9469.68
$ws.Range("F230").Value = 66
$ws.Range("G230").Value = 9469.68

# Row 233: F127->126, G6050.28->6002.64
$ws.Range("F233").Value = 126
$ws.Range("G233").Value = 6002.64

# Row 240: F6->5, G3129.84->2608.2
$ws.Range("F240").Value = 5
$ws.Range("G240").Value = 2608.2

# Row 247: F157->156, G16313.87->16209.96
$ws.Range("F247").Value = 156
$ws.Range("G247").Value = 16209.96

# Row 249: F146->145, G20121.72->19983.9
$ws.Range("F249").Value = 145
$ws.Range("G249").Value = 19983.9

# Row 255: F608->604, G104168.64->103483.32
$ws.Range("F255").Value = 604
$ws.Range("G255").Value = 103483.32

# Row 260: B209722.74->207938
$ws.Range("B260").Value = 207938

# Row 270: F42->37, G1354.08->1192.88
$ws.Range("F270").Value = 37
$ws.Range("G270").Value = 1192.88

# Row 274: F4->1, G1284.04->321.01
$ws.Range("F274").Value = 1
$ws.Range("G274").Value = 321.01

# Row 275: B7668.28->6544.05
$ws.Range("B275").Value = 6544.05

# Row 277: F5->4, G106.25->85
$ws.Range("F277").Value = 4
$ws.Range("G277").Value = 85

# Row 278: F18->16, G2471.76->2197.12
$ws.Range("F278").Value = 16
$ws.Range("G278").Value = 2197.12

# Row 280: F146->145, G24694.44->24525.3
$ws.Range("F280").Value = 145
$ws.Range("G280").Value = 24525.3

# Row 282: F11->10, G590.7->537
$ws.Range("F282").Value = 10
$ws.Range("G282").Value = 537

# Row 283: F47->44, G16049.09->15024.68
$ws.Range("F283").Value = 44
$ws.Range("G283").Value = 15024.68

# Row 285: F15->13, G418.95->363.09
$ws.Range("F285").Value = 13
$ws.Range("G285").Value = 363.09

# Row 288: F53->52, G4928.47->4835.48
$ws.Range("F288").Value = 52
$ws.Range("G288").Value = 4835.48

# Row 291: F123->122, G5290.23->5247.22
$ws.Range("F291").Value = 122
$ws.Range("G291").Value = 5247.22

# Row 294: F49->48, G3496.64->3425.28
$ws.Range("F294").Value = 48
$ws.Range("G294").Value = 3425.28

# Row 296: F93->91, G1971.6->1929.2
$ws.Range("F296").Value = 91
$ws.Range("G296").Value = 1929.2

# Row 299: F278->276, G40209.92->39920.64
$ws.Range("F299").Value = 276
$ws.Range("G299").Value = 39920.64

# Row 302: F78->76, G16449.42->16027.64
$ws.Range("F302").Value = 76
$ws.Range("G302").Value = 16027.64

# Row 304: B196976.59->194416.77
$ws.Range("B304").Value = 194416.77

# Row 306: F73->72, G1541.03->1519.92
$ws.Range("F306").Value = 72
$ws.Range("G306").Value = 1519.92

# Row 309: B1962.8->1941.69
$ws.Range("B309").Value = 1941.69

# Row 328: F59->55, G2195.39->2046.55
$ws.Range("F328").Value = 55
$ws.Range("G328").Value = 2046.55

# Row 330: B32092.07->31943.23
$ws.Range("B330").Value = 31943.23

# Row 338: F83->82, G1967.1->1943.4
$ws.Range("F338").Value = 82
$ws.Range("G338").Value = 1943.4

# Row 345: F79->78, G4851.39->4789.98
$ws.Range("F345").Value = 78
$ws.Range("G345").Value = 4789.98

# Row 346: B28649.43->28564.32
$ws.Range("B346").Value = 28564.32

# Row 364: B53602->65068, E15.69->13.97, F-231->63, G-3037.65->828.45
$ws.Range("B364").Value = 65068
$ws.Range("E364").Value = 13.97
$ws.Range("F364").Value = 63
$ws.Range("G364").Value = 828.45

# Row 365: B65068->53602, E13.97->15.69, F63->-231, G828.45->-3037.65
$ws.Range("B365").Value = 53602
$ws.Range("E365").Value = 15.69
$ws.Range("F365").Value = -231
$ws.Range("G365").Value = -3037.65

# Row 372: B45706->64922, E23.58->20.98, F-202->67, G-3985.46->1321.91
$ws.Range("B372").Value = 64922
$ws.Range("E372").Value = 20.98
$ws.Range("F372").Value = 67
$ws.Range("G372").Value = 1321.91

# Row 373: B64922->45706, E20.98->23.58, F67->-202, G1321.91->-3985.46
$ws.Range("B373").Value = 45706
$ws.Range("E373").Value = 23.58
$ws.Range("F373").Value = -202
$ws.Range("G373").Value = -3985.46

# Row 392: F14->12, G487.34->417.72
$ws.Range("F392").Value = 12
$ws.Range("G392").Value = 417.72

# Row 395: B703.76->634.14
$ws.Range("B395").Value = 634.14

# Row 409: F11->10, G6423.45->5839.5
$ws.Range("F409").Value = 10
$ws.Range("G409").Value = 5839.5

# Row 411: B9979.110000000001->9395.16
$ws.Range("B411").Value = 9395.16

# Row 430: F8->4, G103.12->51.56
$ws.Range("F430").Value = 4
$ws.Range("G430").Value = 51.56

# Row 433: F2->1, G52.76->26.38
$ws.Range("F433").Value = 1
$ws.Range("G433").Value = 26.38

# Row 435: B1593.98->1516.04
$ws.Range("B435").Value = 1516.04

# Row 442: B64810->53319, E291.22->310.64, F5->-6, G1369.6->-1643.52
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52

# Row 443: B53319->64810, E310.64->291.22, F-6->5, G-1643.52->1369.6
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 5
$ws.Range("G443").Value = 1369.6

# Row 447: F3->0, G62.76->0
$ws.Range("F447").Value = 0
$ws.Range("G447").Value = 0

# Row 449: F2->0, G58.06->0
$ws.Range("F449").Value = 0
$ws.Range("G449").Value = 0

# Row 455: F52->51, G3307.72->3244.11
$ws.Range("F455").Value = 51
$ws.Range("G455").Value = 3244.11

# Row 460: B15156.98->14972.55
$ws.Range("B460").Value = 14972.55

# Row 463: B60025->64833, E37.22->34.9, F-98->95, G-3217.34->3118.85
$ws.Range("B463").Value = 64833
$ws.Range("E463").Value = 34.9
$ws.Range("F463").Value = 95
$ws.Range("G463").Value = 3118.85

# Row 464: B64833->60025, E34.9->37.22, F95->-98, G3118.85->-3217.34
$ws.Range("B464").Value = 60025
$ws.Range("E464").Value = 37.22
$ws.Range("F464").Value = -98
$ws.Range("G464").Value = -3217.34

# Row 477: F17->16, G770.78->725.4400000000001
$ws.Range("F477").Value = 16
$ws.Range("G477").Value = 725.4400000000001

# Row 478: B770.78->725.4400000000001
$ws.Range("B478").Value = 725.4400000000001

# Row 485: F29->28, G5088.63->4913.16
$ws.Range("F485").Value = 28
$ws.Range("G485").Value = 4913.16

# Row 488: B33085.9->32910.43
$ws.Range("B488").Value = 32910.43

# Row 509: F249->245, G20014.62->19693.1
$ws.Range("F509").Value = 245
$ws.Range("G509").Value = 19693.1

# Row 510: B26251.02->25929.5
$ws.Range("B510").Value = 25929.5

# Row 550: F6->5, G489.36->407.8
$ws.Range("F550").Value = 5
$ws.Range("G550").Value = 407.8

# Row 552: F28->25, G2850.12->2544.75
$ws.Range("F552").Value = 25
$ws.Range("G552").Value = 2544.75

# Row 560: B9076.77->8689.84
$ws.Range("B560").Value = 8689.84

# Row 575: F5->4, G165.3->132.24
$ws.Range("F575").Value = 4
$ws.Range("G575").Value = 132.24

# Row 577: F82->79, G3525.18->3396.21
$ws.Range("F577").Value = 79
$ws.Range("G577").Value = 3396.21

# Row 578: F101->100, G5038.89->4989
$ws.Range("F578").Value = 100
$ws.Range("G578").Value = 4989

# Row 582: F58->55, G3305.42->3134.45
$ws.Range("F582").Value = 55
$ws.Range("G582").Value = 3134.45

# Row 583: B28970.28->28587.39
$ws.Range("B583").Value = 28587.39

# Row 599: F2065->2037, G336822.15->332255.07
$ws.Range("F599").Value = 2037
$ws.Range("G599").Value = 332255.07

# Row 601: F469->465, G132666.03->131534.55
$ws.Range("F601").Value = 465
$ws.Range("G601").Value = 131534.55

# Row 602: F352->351, G50916.8->50772.15
$ws.Range("F602").Value = 351
$ws.Range("G602").Value = 50772.15

# Row 606: B521253.03->515409.82
$ws.Range("B606").Value = 515409.82

# Row 613: F150->147, G23874->23396.52
$ws.Range("F613").Value = 147
$ws.Range("G613").Value = 23396.52

# Row 618: B46948.53->46471.05
$ws.Range("B618").Value = 46471.05

# Row 619: B1990545.4->1974101.71
$ws.Range("B619").Value = 1974101.71

# Row 620: B1990545.4->1974101.71
$ws.Range("B620").Value = 1974101.71
